# Updates crypto price (column D) and volume change (column E) values
# on the active worksheet, matching the refreshed data pulled on
# Wed Jan 11 15:57:09 UTC 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2 = @{ D = "277.31"; E = "0.54%" }
    3 = @{ D = "27.29"; E = "0.18%" }
    4 = @{ D = "4.848"; E = "1.92%" }
    5 = @{ D = "0.06327"; E = "0.04%" }
    6 = @{ D = "7.031"; E = "1.29%" }
    7 = @{ D = "1.356"; E = "0.04%" }
    8 = @{ D = "0.8912"; E = "1.46%" }
    9 = @{ D = "0.1534"; E = "1.07%" }
    10 = @{ D = "0.05233"; E = "4.30%" }
    11 = @{ D = "0.07417"; E = "-0.98%" }
    12 = @{ D = "0.02899"; E = "-0.01%" }
    13 = @{ D = "0.08943" }
    14 = @{ D = "0.001577"; E = "0.35%" }
    15 = @{ D = "0.0006362"; E = "0.17%" }
    16 = @{ D = "0.006076"; E = "5.17%" }
    17 = @{ D = "3.468"; E = "0.53%" }
    18 = @{ D = "3.304"; E = "0.03%" }
    19 = @{ D = "2.246"; E = "-1.69%" }
    21 = @{ E = "0.71%" }
    22 = @{ D = "3.927"; E = "0.61%" }
    23 = @{ D = "0.1538"; E = "11.43%" }
    24 = @{ D = "0.04402"; E = "0.04%" }
    25 = @{ D = "0.001180"; E = "0.53%" }
    26 = @{ D = "0.004248"; E = "10.59%" }
    28 = @{ D = "0.0001183"; E = "-1.47%" }
    29 = @{ D = "0.0001650"; E = "-14.76%" }
    40 = @{ D = "0.04025"; E = "-2.18%" }
    41 = @{ D = "0.006829"; E = "0.04%" }
    42 = @{ D = "0.1412"; E = "20.54%" }
    43 = @{ D = "0.002075"; E = "6.92%" }
    44 = @{ D = "0.01105"; E = "-3.71%" }
    45 = @{ D = "0.00005357"; E = "2.24%" }
    46 = @{ D = "1.628"; E = "9.30%" }
    47 = @{ D = "0.01852"; E = "-7.33%" }
}

foreach ($rowKey in $changes.Keys) {
    $rowChanges = $changes[$rowKey]
    foreach ($col in $rowChanges.Keys) {
        $addr = "$col$rowKey"
        $cell = $ws.Range($addr)
        # Force text storage so values like "277.31" and "0.54%"
        # stay as literal strings instead of being parsed into
        # numbers / percentages by Excel.
        $cell.NumberFormat = "@"
        $cell.Value = $rowChanges[$col]
        # Restore the default (unstyled) cell formatting so no stray
        # number-format styling is left behind on the cell.
        $cell.Style = "Normal"
    }
}

Write-Host "Updated $($changes.Count) rows of crypto price/volume data"
